# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" detail sheet (cloned from the "2022-Q2" sheet so it
# keeps identical formatting) right after the "总计" summary sheet, fills it
# in with the new quarter's numbers, and adds a matching new row at the top
# of the "总计" summary table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet by cloning "2022-Q2" (so borders,
#    fonts, column layout, etc. all match the other quarterly sheets) and
#    position it immediately before "2022-Q2".
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($template, $null)
$newSheet = $wb.ActiveSheet
$newSheet.Name = "2022-Q3"

# Overwrite the cloned data row with the new quarter's figures. The four
# numeric-looking columns (D/E/F/G) are stored as text in the source data,
# so force a leading quote to keep Excel from re-typing them as numbers
# (which would strip the trailing zeros, e.g. "2.20" -> 2.2).
$newSheet.Cells(2, 4).Value = "'0.41"
$newSheet.Cells(2, 5).Value = "'91.47"
$newSheet.Cells(2, 6).Value = "'2.20"
$newSheet.Cells(2, 7).Value = "'0.0090"
$newSheet.Cells(2, 8).Value = 7

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: push the existing rows down by one
#    and insert the new "2022-Q3" row at the top of the table.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Remember the existing rows (row 2 = newest … row 4 = oldest) before they
# get overwritten.
$oldLabel2 = $summary.Cells(2, 2).Value()
$oldCount2 = $summary.Cells(2, 3).Value()
$oldValue2 = $summary.Cells(2, 4).Value()

$oldLabel3 = $summary.Cells(3, 2).Value()
$oldCount3 = $summary.Cells(3, 3).Value()
$oldValue3 = $summary.Cells(3, 4).Value()

$oldLabel4 = $summary.Cells(4, 2).Value()
$oldCount4 = $summary.Cells(4, 3).Value()
$oldValue4 = $summary.Cells(4, 4).Value()

# Clone the formatting of the last existing data row (row 4, column A) down
# onto the new row 5 so the index column keeps its style.
$summary.Cells(4, 1).Copy($summary.Cells(5, 1))

# Row 2 becomes the brand-new "2022-Q3" entry.
$summary.Cells(2, 1).Value = 0
$summary.Cells(2, 2).Value = "2022-Q3"
$summary.Cells(2, 3).Value = 1
$summary.Cells(2, 4).Value = 0.01

# Rows 3 and 4 shift down to what used to be rows 2 and 3.
$summary.Cells(3, 1).Value = 1
$summary.Cells(3, 2).Value = $oldLabel2
$summary.Cells(3, 3).Value = $oldCount2
$summary.Cells(3, 4).Value = $oldValue2

$summary.Cells(4, 1).Value = 2
$summary.Cells(4, 2).Value = $oldLabel3
$summary.Cells(4, 3).Value = $oldCount3
$summary.Cells(4, 4).Value = $oldValue3

# Row 5 is new, holding what used to be row 4's data.
$summary.Cells(5, 1).Value = 3
$summary.Cells(5, 2).Value = $oldLabel4
$summary.Cells(5, 3).Value = $oldCount4
$summary.Cells(5, 4).Value = $oldValue4
